{"js": "// Update calibration model in main analysis\n//\n// Table 1 (\"Performance metrics based on the leave-one-sensor-out\n// approach\") has two \"Calibrated\" rows (DustTrak and SMPS + APS). Their\n// Mean Relative Error (%), Mean Error, LoA Low and LoA High figures were\n// refreshed with new calibration results. Every old value is a unique\n// string in the document, so searching the body for each exact value and\n// replacing it in place is safe and leaves all other cells/formatting\n// untouched.\n\nconst body = context.document.body;\n\nasync function replaceValue(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// DustTrak / Calibrated row\nawait replaceValue(\"4.42\", \"6.12\");\nawait replaceValue(\"-1.16\", \"-1.11\");\nawait replaceValue(\"-34.34\", \"-34.08\");\nawait replaceValue(\"32.02\", \"31.86\");\n\n// SMPS + APS / Calibrated row\nawait replaceValue(\"72.9\", \"72.86\");\nawait replaceValue(\"-6.19\", \"-6.22\");\nawait replaceValue(\"-195.14\", \"-195.19\");\nawait replaceValue(\"182.75\", \"182.74\");\n", "ps1": "# Update calibration model in main analysis\n# Table 1 (Performance metrics): for the two \"Calibrated\" rows, the Mean\n# Relative Error (%), Mean Error, LoA Low and LoA High figures were\n# refreshed. Each old value is unique in the document, so a plain\n# whole-word Find/Replace for each pair is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Value($oldText, $newText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $null = $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\n# DustTrak / Calibrated row\nReplace-Value \"4.42\" \"6.12\"\nReplace-Value \"-1.16\" \"-1.11\"\nReplace-Value \"-34.34\" \"-34.08\"\nReplace-Value \"32.02\" \"31.86\"\n\n# SMPS + APS / Calibrated row\nReplace-Value \"72.9\" \"72.86\"\nReplace-Value \"-6.19\" \"-6.22\"\nReplace-Value \"-195.14\" \"-195.19\"\nReplace-Value \"182.75\" \"182.74\"\n"}
